# Add common skill "Hymn of Patience": introduces a new "Mine" column that
# tracks the author's own custom-character deck alongside Silent/Ironclad/Defect
# in all three distribution tables on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Cost-distribution table (Table2, A6:D13) -> add "Mine" column (E)
# ---------------------------------------------------------------------
$costTable = $ws.ListObjects.Item("Table2")
$costTable.ListColumns.Add() | Out-Null

$ws.Range("E6").Value = "Mine"
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 6
$ws.Range("E9").Value = 18
$ws.Range("E10").Value = 10
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 0
$ws.Range("E13").Formula = "=SUM(Table2[Mine])"

$mineCostRange = $ws.Range("E6:E13")
$mineCostRange.HorizontalAlignment = -4108
$mineCostRange.VerticalAlignment = -4108
$mineCostRange.WrapText = $true

# ---------------------------------------------------------------------
# 2) Rarity-distribution block (Table24, F6:I11) -> add "Mine" column (J)
#    (kept as plain worksheet cells mirroring the table, same as the
#    existing Silent/Ironclad/Defect columns' neighbour layout)
# ---------------------------------------------------------------------
$ws.Range("J6").Value = "Mine"
$ws.Range("J7").Value = 4
$ws.Range("J8").Value = 10
$ws.Range("J9").Value = 16
$ws.Range("J10").Value = 10
$ws.Range("J11").Formula = "=SUM(J7:J10)"

# ---------------------------------------------------------------------
# 3) Type-distribution block (Table242, F16:I20) -> add "Mine" column (J)
# ---------------------------------------------------------------------
$ws.Range("J16").Value = "Mine"
$ws.Range("J17").Value = 17
$ws.Range("J18").Value = 19
$ws.Range("J19").Value = 4
$ws.Range("J20").Formula = "=SUM(J17:J19)"

$mineJRange = $ws.Range("J6:J20")
$mineJRange.HorizontalAlignment = -4108
$mineJRange.VerticalAlignment = -4108
$mineJRange.WrapText = $true

# Match the column width of the new "Mine" column (J) to its neighbours.
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# Final selection, matching the state the workbook was left in.
$ws.Range("J9").Select()
